$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ebf1d1d5be3be99f5437917afe2f55abdc20e83/e2e/ad2609dc-a353-4226-b0f0-a908752174ce.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb648a8ee87d6834c36c8bb14172ac899c46ca7a/e2e/ad2609dc-a353-4226-b0f0-a908752174ce.md."

# Overview sheet, row 3 ("ad2609dc-a353-4226-b0f0-a908752174ce.md")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = "2016-09-06 10:36:29"

# zh-cn sheet, row 3
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-09-06 10:36:16"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de sheet, row 3
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-09-06 10:36:29"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
